$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct part name in BOM: 2N7002 -> 2N7000 (row 8, column A)
$ws.Range("A8").Value = "2N7000"
